# CIERRE 3 SEPT 22
# Roll the payroll receipt ("recibos") sheet forward from week 34
# (Aug 22-28, 2022) to week 35 (Aug 29 - Sep 04, 2022).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("recibos")   # same sheet as $wb.ActiveSheet

# --- Week label (B9) -------------------------------------------------
# H9 (=B9), B27 (=B9), H27 (=B27) and B43 (=H27) all recompute from this.
$ws.Range("B9").Value = "SEMANA  35  DEL    29      Al   04   DE  SEPTIEMBRE          2022"

# --- Extra / bonus paid this period (K21) -----------------------------
# K24 (=SUM(K21:K23)) recomputes automatically.
$ws.Range("K21").Value = 560

# --- Horas extra: days worked & amount (D38 / E38) --------------------
# E41 (=SUM(E38:E40)) recomputes automatically.
$ws.Range("D38").Value = 6
$ws.Range("E38").Value = 2500

# --- Scroll position / selection left by the editor --------------------
$ws.Activate() | Out-Null
$excel.ActiveWindow.ScrollRow = 22
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("E39").Select() | Out-Null
